$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$table = $ws.ListObjects.Item("Table1")

$newDate = Get-Date -Year 2024 -Month 11 -Day 10 -Hour 0 -Minute 0 -Second 0

# New row 1: Spez 1 / Kaffee und Gipfeli
$row1 = $table.ListRows.Add()
$row1.Range.Item(1, 1).Value = $newDate
$row1.Range.Item(1, 2).Value = "Spez 1"
$row1.Range.Item(1, 3).Value = "Kaffee und Gipfeli"
$row1.Range.Item(1, 4).Value = 5
$row1.Range.Item(1, 5).Value = 9

# New row 2: Spez 2 / Gipfeli
$row2 = $table.ListRows.Add()
$row2.Range.Item(1, 1).Value = $newDate
$row2.Range.Item(1, 2).Value = "Spez 2"
$row2.Range.Item(1, 3).Value = "Gipfeli"
$row2.Range.Item(1, 4).Value = 1
$row2.Range.Item(1, 5).Value = 1

$ws.Range("E29").Select()
